# "added test cases and utils"
# Turns the single-cell Sheet1 (A1 only) into a two-column test-case table
# (TC_Name / TypeValue headers + 5 data rows), matching the target workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A header
$ws.Range("A1").Value = "TC_Name"

# Column B data values (written in this order so the resulting shared-string
# table lines up with the target workbook's string order)
$ws.Range("B3").Value = "gfr6657fr6fg"
$ws.Range("B2").Value = "mobiles"
$ws.Range("B4").Value = "chocolate gift hamper"
$ws.Range("B5").Value = "panda teddy"
$ws.Range("B6").Value = "bag"

# Column A test-case descriptions
$ws.Range("A2").Value = "Verify that the user can able to search a product"
$ws.Range("A3").Value = "Verify user can able to search invalid product"
$ws.Range("A4").Value = "Verify user search chocolate gift hamper product"
$ws.Range("A5").Value = "Verify user search panda teddy product"
$ws.Range("A6").Value = "Verify user search  bag product"

# Column B header
$ws.Range("B1").Value = "TypeValue"

# Widen the two columns to fit the new content
$ws.Columns.Item(1).ColumnWidth = 64.833333333
$ws.Columns.Item(2).ColumnWidth = 52.666666667

# Leave the same cell selected as in the target workbook
$ws.Range("B4").Select()
